$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "ıd"
$ws.Range("B1").Value = "numaralar"
$ws.Range("B1").Select()
